# "Add files via upload" - the Global sheet's credential table is replaced
# with a minimal single-row A/B header, and the leftover demo rows/styles
# are removed.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Global" sheet
$ws2 = $wb.Worksheets.Item(2)   # "Action2" sheet (originally active)

# Write the new header values. B1 is set before A1 so the shared-string
# table ends up ordered ["B", "A"], matching the target workbook.
$ws1.Range("B1").Value = "B"
$ws1.Range("A1").Value = "A"

# Remove the old demo rows (username/aidemo/user2/user3 + passwords),
# shrinking the used range down to just row 1.
$ws1.Rows("2:4").Delete() | Out-Null

# Update the remembered selection on the Global sheet.
$ws1.Range("C2").Select() | Out-Null

# Selecting on ws1 makes it the active sheet as a side effect; restore the
# original active sheet (Action2) so its tab/activation state is unchanged.
$ws2.Activate() | Out-Null
